$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns with the latest
# scraped cryptocurrency quotes. Values that would otherwise be
# auto-parsed by Excel as numbers (losing significant trailing zeros,
# e.g. "1.140" -> 1.14) are entered with a leading apostrophe so Excel
# keeps them as literal text, matching the source data feed.
$ws.Range("D2").Value = '28.089.97'
$ws.Range("E2").Value = '  -1.81%  '
$ws.Range("D3").Value = '1.835.06'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''325.45'
$ws.Range("E5").Value = '  -3.21%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '''0.4616'
$ws.Range("E7").Value = '  -1.03%  '
$ws.Range("D8").Value = '''0.3863'
$ws.Range("E8").Value = '  -0.95%  '
$ws.Range("D9").Value = '''0.07847'
$ws.Range("E9").Value = '  -0.83%  '
$ws.Range("D10").Value = '''0.9617'
$ws.Range("E10").Value = '  -1.80%  '
$ws.Range("D11").Value = '''22.02'
$ws.Range("E11").Value = '  -1.20%  '
$ws.Range("D12").Value = '1.808.15'
$ws.Range("E12").Value = '  -1.16%  '
$ws.Range("D13").Value = '''5.677'
$ws.Range("E13").Value = '  -2.52%  '
$ws.Range("D14").Value = '''6.900'
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("D15").Value = '''0.06844'
$ws.Range("E15").Value = '  -1.11%  '
$ws.Range("D16").Value = '''88.38'
$ws.Range("E16").Value = '  +0.88%  '
$ws.Range("D17").Value = '''1.000'
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("D18").Value = '''0.000009934'
$ws.Range("E18").Value = '  -0.92%  '
$ws.Range("D19").Value = '''16.66'
$ws.Range("E19").Value = '  -2.32%  '
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("D21").Value = '28.097.49'
$ws.Range("E21").Value = '  -1.80%  '
$ws.Range("D22").Value = '''5.292'
$ws.Range("E22").Value = '  -1.83%  '
$ws.Range("D23").Value = '''11.02'
$ws.Range("E23").Value = '  -2.29%  '
$ws.Range("D24").Value = '''2.087'
$ws.Range("E24").Value = '  -3.93%  '
$ws.Range("D25").Value = '2.074.31'
$ws.Range("E25").Value = '  +1.00%  '
$ws.Range("D26").Value = '''154.66'
$ws.Range("E26").Value = '  +0.85%  '
$ws.Range("D27").Value = '''19.13'
$ws.Range("E27").Value = '  -1.40%  '
$ws.Range("D28").Value = '''5.732'
$ws.Range("E28").Value = '  -5.33%  '
$ws.Range("D29").Value = '''1.978'
$ws.Range("E29").Value = '  -2.34%  '
$ws.Range("D30").Value = '''119.69'
$ws.Range("E30").Value = '  +1.96%  '
$ws.Range("D31").Value = '''0.9405'
$ws.Range("E31").Value = '  -3.07%  '
$ws.Range("E32").Value = '  -0.91%  '
$ws.Range("D33").Value = '''5.271'
$ws.Range("E33").Value = '  -1.55%  '
$ws.Range("E34").Value = '  -1.84%  '
$ws.Range("E35").Value = '  -4.42%  '
$ws.Range("D36").Value = '''0.05830'
$ws.Range("E36").Value = '  -5.29%  '
$ws.Range("D37").Value = '''0.02113'
$ws.Range("E37").Value = '  -3.94%  '
$ws.Range("D38").Value = '''1.140'
$ws.Range("E38").Value = '  -2.52%  '
$ws.Range("D39").Value = '''7.736'
$ws.Range("E39").Value = '  +0.77%  '
$ws.Range("D40").Value = '''0.5603'
$ws.Range("E40").Value = '  -1.86%  '
$ws.Range("D41").Value = '''9.909'
$ws.Range("E41").Value = '  -2.33%  '
$ws.Range("D43").Value = '''0.07325'
$ws.Range("E43").Value = '  +3.18%  '
$ws.Range("D44").Value = '''11.63'
$ws.Range("E44").Value = '  -0.94%  '
$ws.Range("D45").Value = '''0.5268'
$ws.Range("E45").Value = '  -2.10%  '
$ws.Range("D46").Value = '''2.128'
$ws.Range("E46").Value = '  -11.98%  '
$ws.Range("D47").Value = '''1.132'
$ws.Range("E47").Value = '  -9.40%  '
$ws.Range("D48").Value = '''1.837'
$ws.Range("E48").Value = '  -3.53%  '
$ws.Range("D49").Value = '''113.55'
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("E50").Value = '  -0.11%  '
$ws.Range("E51").Value = '  +0.10%  '
